$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 16-17 (pushes existing rows 16.. down to 18..,
# extending the used range from A1:T126 to A1:T128).
$ws.Rows("16:17").Insert()

# New row 16: Packham's Triumph / Primera, same market/variety as the row
# that follows, new sampling date (2021-09-23) and a smaller Volumen (200).
$ws.Cells.Item(16, 1).Value = 4
$ws.Cells.Item(16, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(16, 3).Value = "Los Lagos"
$ws.Cells.Item(16, 4).Value = 44462
$ws.Cells.Item(16, 5).Value = 10
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100104
$ws.Cells.Item(16, 8).Value = "Frutos de pepita"
$ws.Cells.Item(16, 9).Value = 100104005
$ws.Cells.Item(16, 10).Value = "Pera"
$ws.Cells.Item(16, 11).Value = "Packham's Triumph"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 16000
$ws.Cells.Item(16, 15).Value = 16000
$ws.Cells.Item(16, 16).Value = 16000
$ws.Cells.Item(16, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 1067
$ws.Cells.Item(16, 20).Value = 15

# New row 17: Packham's Triumph / Segunda, same new sampling date, smaller
# Volumen (100).
$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value = "Los Lagos"
$ws.Cells.Item(17, 4).Value = 44462
$ws.Cells.Item(17, 5).Value = 10
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100104
$ws.Cells.Item(17, 8).Value = "Frutos de pepita"
$ws.Cells.Item(17, 9).Value = 100104005
$ws.Cells.Item(17, 10).Value = "Pera"
$ws.Cells.Item(17, 11).Value = "Packham's Triumph"
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 12000
$ws.Cells.Item(17, 15).Value = 12000
$ws.Cells.Item(17, 16).Value = 12000
$ws.Cells.Item(17, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(17, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 19).Value = 800
$ws.Cells.Item(17, 20).Value = 15
